$d = $word.ActiveDocument

# Locate the paragraph that contains the "LOQ4087..." requirement line,
# and the paragraph that contains the copyright notice; remove everything
# from right after the former up to (and including) the end of the latter.
# That deletes the blank paragraph, the "Ver no Jupiter..." paragraph, and
# the "(c) 2020 ..." paragraph, leaving the requirement line followed
# directly by the (already existing) trailing blank paragraph.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOQ4087*") {
        $startPara = $p
    }
    if ($t -like "*Powered by Jekyll*") {
        $endPara = $p
        break
    }
}

$deleteRange = $d.Range($startPara.Range.End, $endPara.Range.End)
$deleteRange.Delete()
